# Actualización automática 2025-06-17 17:30:08
$wb = $excel.ActiveWorkbook

# Target sheet: "CUMPLIMIENTO MENSUAL"
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Narrow column F slightly (stored width 26 -> 25).
# Note: the ColumnWidth property adds a constant padding offset (5/6)
# before the value is persisted as the sheet's stored "width" attribute,
# so we subtract that offset here to land exactly on a stored width of 25.
$ws.Columns.Item(6).ColumnWidth = (25 - 5/6)

# Row 2 (OTROS): update VENTA (D) and recompute POR CUMPLIR (E)
$ws.Range("D2").Value = 2942.59
$ws.Range("E2").Value = -2942.59

# Row 3 (PORCELANATO): update VENTA (D), POR CUMPLIR (E) and CUMPLIMIENTO (F)
$ws.Range("D3").Value = 248.48
$ws.Range("E3").Value = 17251.52
$ws.Range("F3").Value = 0.01419885714285714

# Row 4 (TOTAL): update VENTA (D), POR CUMPLIR (E) and CUMPLIMIENTO (F)
$ws.Range("D4").Value = 3191.07
$ws.Range("E4").Value = 14308.93
$ws.Range("F4").Value = 0.1823468571428571
